# Letter Model inside the CheckIn Pages
#
# This script turns:
#   Donation page
#   Anchor tag
#   Comment for each letter
#   Picture carousel
#   Contact Message
#   Admin account [_GoBack bookmark]
# into:
#   Donation page [highlighted yellow]
#   Comment for each letter [highlighted yellow]
#   Picture carousel [_GoBack bookmark]
#   Admin account

$d = $word.ActiveDocument

# Locate the relevant paragraphs by their (unique) text so the script
# does not depend on brittle, hard-coded paragraph indices.
function Get-ParaByText($doc, $text) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Trim() -eq $text) {
            return $p
        }
    }
    return $null
}

# 1. Highlight "Donation page" in yellow. Using Range.Font.HighlightColorIndex
#    (rather than Range.HighlightColorIndex) also stamps the paragraph mark's
#    run properties (w:pPr/w:rPr), matching the authored formatting.
$pDonation = Get-ParaByText $d "Donation page"
$pDonation.Range.Font.HighlightColorIndex = 7

# 2. Remove the "Anchor tag" paragraph entirely.
$pAnchor = Get-ParaByText $d "Anchor tag"
$pAnchor.Range.Delete()

# 3. Highlight "Comment for each letter" in yellow (run only, no pPr mark).
$pComment = Get-ParaByText $d "Comment for each letter"
$pComment.Range.HighlightColorIndex = 7

# 4. Move the hidden "_GoBack" bookmark so that it now sits at the end of the
#    "Picture carousel" paragraph (after its run) instead of on "Admin
#    account". Re-using the name "_GoBack" for Bookmarks.Add automatically
#    relocates (rather than duplicates) the existing bookmark.
#    A zero-length range placed directly before a paragraph mark is not
#    reliably positioned by this host, so a placeholder character is
#    inserted, bookmarked, and then removed again -- leaving the bookmark
#    collapsed at the desired location.
$pPicture = Get-ParaByText $d "Picture carousel"
$insertPos = $pPicture.Range.End - 1
$placeholder = $d.Range($insertPos, $insertPos)
$placeholder.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $placeholder)
$bm = $d.Bookmarks("_GoBack")
$bm.Range.Text = ""

# 5. Remove the "Contact Message" paragraph entirely.
$pContact = Get-ParaByText $d "Contact Message"
$pContact.Range.Delete()
